$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.796.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'1.758.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'237.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.5064"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'0.2663"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.16%  "
$ws.Range("D10").Value = "'0.06206"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("D11").Value = "'1.752.60"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'0.06941"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").Value = "'15.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +9.85%  "
$ws.Range("D14").Value = "'0.6065"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'4.468"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("D16").Value = "'77.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "'25.839.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'0.000006828"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.08%  "
$ws.Range("D21").Value = "'11.63"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.22%  "
$ws.Range("D22").Value = "'1.980.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'4.068"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.75%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.198"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'8.167"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").Value = "'137.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").Value = "'1.457"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "'15.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("D30").Value = "'102.67"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("D31").Value = "'0.08233"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").Value = "'3.687"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("E33").Value = "  +8.43%  "
$ws.Range("D34").Value = "'0.04374"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("D35").Value = "'0.9993"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").Value = "'2.654"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'0.6077"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'2.728"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'1.937"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").Value = "'0.01545"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.94%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'103.19"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").Value = "'0.3834"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "'0.7385"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").Value = "'4.906"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("D47").Value = "'0.05491"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.53%  "
$ws.Range("D48").Value = "'0.1080"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.08%  "
$ws.Range("D49").Value = "'5.943"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'29.93"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "'7.607"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.92%  "
